$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update launch times (column B, rows 2-4); leading apostrophe keeps them as
# quote-prefixed text entries (same as the original cells).
$ws.Range("B2").Formula = "'11 Apr 2018 16:00:00.000'"
$ws.Range("B3").Formula = "'11 Apr 2018 18:00:00.000'"
$ws.Range("B4").Formula = "'11 Apr 2018 19:00:00.000'"

# Update latitude/longitude values
$ws.Range("C2").Value = 33
$ws.Range("D2").Value = -104

$ws.Range("C3").Value = 33.5
$ws.Range("D3").Value = -104.5

$ws.Range("C4").Value = 34
$ws.Range("D4").Value = -103

# Update the selection on the sheet
$ws.Range("D4").Select()
